$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to store values as text so that numeric-looking
# strings (e.g. "320.56") are not silently converted to numbers, matching the
# original inline-string cell contents. The temporary text format is reverted
# to the default "Normal" style afterwards so no stray formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "47.186.73"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "2.484.49"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "320.56"
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").Value = "107.48"
$ws.Range("E6").Value = "  +1.32%  "

$ws.Range("D7").Value = "0.520"
$ws.Range("E7").Value = "  -0.66%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -1.55%  "

$ws.Range("D10").Value = "38.45"
$ws.Range("E10").Value = "  +4.95%  "

$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("D13").Value = "18.18"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("D14").Value = "7.10"
$ws.Range("E14").Value = "  -0.98%  "

$ws.Range("D15").Value = "2.876.35"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").Value = "2.485.13"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("D17").Value = "0.844"
$ws.Range("E17").Value = "  -0.43%  "

$ws.Range("D18").Value = "47.121.90"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  -0.32%  "

$ws.Range("D20").Value = "6.60"
$ws.Range("E20").Value = "  +1.26%  "

$ws.Range("E21").Value = "  -1.12%  "

$ws.Range("D22").Value = "2.70"
$ws.Range("E22").Value = "  +12.72%  "

$ws.Range("D23").Value = "70.26"
$ws.Range("E23").Value = "  -0.99%  "

$ws.Range("D24").Value = "245.08"
$ws.Range("E24").Value = "  -3.11%  "

$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "25.63"
$ws.Range("E27").Value = "  -3.11%  "

$ws.Range("E28").Value = "  +3.21%  "

$ws.Range("D29").Value = "9.98"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("D30").Value = "34.35"
$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("D31").Value = "0.135"
$ws.Range("E31").Value = "  -1.59%  "

$ws.Range("D32").Value = "49.52"

$ws.Range("D33").Value = "20.15"
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("D34").Value = "5.33"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "0.0778"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("D37").Value = "1.95"
$ws.Range("E37").Value = "  +0.62%  "

$ws.Range("D38").Value = "4.60"
$ws.Range("E38").Value = "  -0.96%  "

$ws.Range("E39").Value = "  -1.11%  "

$ws.Range("E40").Value = "  -0.74%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "22.32"
$ws.Range("E41").Value = "  +6.11%  "

$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  -0.65%  "

$ws.Range("D43").Value = "119.06"
$ws.Range("E43").Value = "  -4.35%  "

$ws.Range("D44").Value = "0.0294"
$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("D45").Value = "1.981.39"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  -1.10%  "

$ws.Range("D47").Value = "1.98"
$ws.Range("E47").Value = "  -6.30%  "

$ws.Range("D48").Value = "9.03"
$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").Value = "1.75"
$ws.Range("E49").Value = "  -2.67%  "

$ws.Range("D50").Value = "5.09"
$ws.Range("E50").Value = "  -5.88%  "

$ws.Range("D51").Value = "56.42"
$ws.Range("E51").Value = "  +2.91%  "

# Restore the default cell style on the Price column now that the values are set,
# so only the cell contents changed (no lingering number-format override).
$ws.Range("D2:D51").Style = "Normal"
